$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1637.6666
$ws.Range("J17").Value = 1637.6666
$ws.Range("L17").Value = 4912.9998
$ws.Range("N17").Value = -5248.9998
$ws.Range("H28").Value = 2173.7273
$ws.Range("I28").Value = 434.55554
$ws.Range("J28").Value = 10000
$ws.Range("K28").Value = 434.55554
$ws.Range("L28").Value = 10000
$ws.Range("M28").Value = 50.44445999999999
$ws.Range("N28").Value = -10970
$ws.Range("H69").Value = 19249.916
$ws.Range("I69").Value = 18000
$ws.Range("J69").Value = 19666.555
$ws.Range("K69").Value = 54000
$ws.Range("L69").Value = 58999.665
$ws.Range("M69").Value = -53126
$ws.Range("N69").Value = -60747.665
$ws.Range("H72").Value = 19249.916
$ws.Range("I72").Value = 18000
$ws.Range("J72").Value = 19666.555
$ws.Range("K72").Value = 162000
$ws.Range("L72").Value = 176998.995
$ws.Range("M72").Value = -157632
$ws.Range("N72").Value = -185734.995
$ws.Range("H76").Value = 10051.556
$ws.Range("J76").Value = 8435.625
$ws.Range("L76").Value = 8435.625
$ws.Range("N76").Value = -9065.625
$ws.Range("H79").Value = 10051.556
$ws.Range("J79").Value = 8435.625
$ws.Range("L79").Value = 8435.625
$ws.Range("N79").Value = -10619.625
$ws.Range("H80").Value = 1293.0667
$ws.Range("I80").Value = 777
$ws.Range("J80").Value = 1480.7273
$ws.Range("K80").Value = 2331
$ws.Range("L80").Value = 4442.1819
$ws.Range("M80").Value = -1333
$ws.Range("N80").Value = -6438.1819
$ws.Range("H82").Value = 4487.625
$ws.Range("I82").Value = 1316.8334
$ws.Range("J82").Value = 14000
$ws.Range("K82").Value = 3950.5002
$ws.Range("L82").Value = 42000
$ws.Range("M82").Value = -3544.5002
$ws.Range("N82").Value = -42812
$ws.Range("H83").Value = 1293.0667
$ws.Range("I83").Value = 777
$ws.Range("J83").Value = 1480.7273
$ws.Range("K83").Value = 6993
$ws.Range("L83").Value = 13326.5457
$ws.Range("M83").Value = -2001
$ws.Range("N83").Value = -23310.5457
$ws.Range("H85").Value = 4487.625
$ws.Range("I85").Value = 1316.8334
$ws.Range("J85").Value = 14000
$ws.Range("K85").Value = 3950.5002
$ws.Range("L85").Value = 42000
$ws.Range("M85").Value = -2546.5002
$ws.Range("N85").Value = -44808
$ws.Range("H116").Value = 5988.5884
$ws.Range("I116").Value = 5329.75
$ws.Range("K116").Value = 5329.75
$ws.Range("M116").Value = -1887.75
$ws.Range("H135").Value = 7357
$ws.Range("I135").Value = 1300.4
$ws.Range("J135").Value = 16009.286
$ws.Range("K135").Value = 11703.6
$ws.Range("L135").Value = 144083.574
$ws.Range("M135").Value = -9168.6
$ws.Range("N135").Value = -149153.574
$ws.Range("H137").Value = 1871.75
$ws.Range("J137").Value = 2774.6667
$ws.Range("L137").Value = 8324.000100000001
$ws.Range("N137").Value = -13424.0001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H94").Value = 144665
$ws.Range("J94").Value = 144665
$ws.Range("L94").Value = 144665
$ws.Range("N94").Value = -146467
$ws.Range("H110").Value = 1721
$ws.Range("I110").Value = 1616.619
$ws.Range("K110").Value = 1616.619
$ws.Range("M110").Value = 428.3810000000001
$ws.Range("H122").Value = 2557.1875
$ws.Range("I122").Value = 2235.923
$ws.Range("K122").Value = 6707.768999999999
$ws.Range("M122").Value = -4257.768999999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3281.65
$ws.Range("I134").Value = 2927.3157
$ws.Range("K134").Value = 8781.947100000001
$ws.Range("M134").Value = -6246.947100000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2301
$ws.Range("I58").Value = 2308.75
$ws.Range("J58").Value = 2270
$ws.Range("K58").Value = 2308.75
$ws.Range("L58").Value = 2270
$ws.Range("M58").Value = -2105.75
$ws.Range("N58").Value = -2676
$ws.Range("H99").Value = 3942.5789
$ws.Range("I99").Value = 3800
$ws.Range("J99").Value = 4070.9
$ws.Range("K99").Value = 3800
$ws.Range("L99").Value = 4070.9
$ws.Range("M99").Value = -2302
$ws.Range("N99").Value = -7066.9
$ws.Range("H105").Value = 1305.5
$ws.Range("I105").Value = 1117.2222
$ws.Range("K105").Value = 1117.2222
$ws.Range("M105").Value = 629.7778000000001
$ws.Range("H122").Value = 3851.8823
$ws.Range("I122").Value = 3904.25
$ws.Range("J122").Value = 3726.2
$ws.Range("K122").Value = 11712.75
$ws.Range("L122").Value = 11178.6
$ws.Range("M122").Value = -9262.75
$ws.Range("N122").Value = -16078.6
$ws.Range("H126").Value = 3942.5789
$ws.Range("I126").Value = 3800
$ws.Range("J126").Value = 4070.9
$ws.Range("K126").Value = 11400
$ws.Range("L126").Value = 12212.7
$ws.Range("M126").Value = -8930
$ws.Range("N126").Value = -17152.7
$ws.Range("H136").Value = 2301
$ws.Range("I136").Value = 2308.75
$ws.Range("J136").Value = 2270
$ws.Range("K136").Value = 6926.25
$ws.Range("L136").Value = 6810
$ws.Range("M136").Value = -4376.25
$ws.Range("N136").Value = -11910

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 116.27273
$ws.Range("I33").Value = 109
$ws.Range("J33").Value = 122.333336
$ws.Range("K33").Value = 654
$ws.Range("L33").Value = 734.000016
$ws.Range("M33").Value = -371
$ws.Range("N33").Value = -1300.000016
$ws.Range("H92").Value = 407.33334
$ws.Range("I92").Value = 398.5
$ws.Range("J92").Value = 411.75
$ws.Range("K92").Value = 1195.5
$ws.Range("L92").Value = 1235.25
$ws.Range("M92").Value = 52.5
$ws.Range("N92").Value = -3731.25
$ws.Range("H100").Value = 14285.571
$ws.Range("J100").Value = 14285.571
$ws.Range("L100").Value = 42856.713
$ws.Range("N100").Value = -44478.713
$ws.Range("H113").Value = 2066.7778
$ws.Range("I113").Value = 1930.2
$ws.Range("J113").Value = 2237.5
$ws.Range("K113").Value = 5790.6
$ws.Range("L113").Value = 6712.5
$ws.Range("M113").Value = -3620.6
$ws.Range("N113").Value = -11052.5
$ws.Range("H131").Value = 1859.1025
$ws.Range("I131").Value = 1199.2667
$ws.Range("J131").Value = 2271.5
$ws.Range("K131").Value = 3597.800099999999
$ws.Range("L131").Value = 6814.5
$ws.Range("M131").Value = 1442.199900000001
$ws.Range("N131").Value = -16894.5
$ws.Range("H132").Value = 2115.9092
$ws.Range("I132").Value = 998
$ws.Range("J132").Value = 2227.7
$ws.Range("K132").Value = 8982
$ws.Range("L132").Value = 20049.3
$ws.Range("M132").Value = -6452
$ws.Range("N132").Value = -25109.3

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 234.625
$ws.Range("I2").Value = 235.4
$ws.Range("K2").Value = 235.4
$ws.Range("M2").Value = -122.4
$ws.Range("H80").Value = 2503.2
$ws.Range("I80").Value = 2503.2
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2503.2
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1505.2
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 2503.2
$ws.Range("I83").Value = 2503.2
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 12516
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -7524
$ws.Range("N83").ClearContents()
$ws.Range("H126").Value = 5759.5
$ws.Range("I126").Value = 6023.4
$ws.Range("J126").Value = 5495.6
$ws.Range("K126").Value = 18070.2
$ws.Range("L126").Value = 16486.8
$ws.Range("M126").Value = -15600.2
$ws.Range("N126").Value = -21426.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 6682.6665
$ws.Range("I55").Value = 7524
$ws.Range("K55").Value = 7524
$ws.Range("M55").Value = -7247
$ws.Range("H107").Value = 3478.8823
$ws.Range("I107").Value = 3224.7144
$ws.Range("K107").Value = 9674.143199999999
$ws.Range("M107").Value = -7754.143199999999
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("N123").ClearContents()
$ws.Range("H126").Value = 3244.923
$ws.Range("I126").Value = 2093.5
$ws.Range("J126").Value = 3454.2727
$ws.Range("K126").Value = 6280.5
$ws.Range("L126").Value = 10362.8181
$ws.Range("M126").Value = -3810.5
$ws.Range("N126").Value = -15302.8181
